# Fixed naive component forecaster bug - Presentation state 11.02.
#
# A newly computed "Q0" error value is inserted at the front of each data
# row (column B). The values that used to live in B:K shift one column to
# the right (B->C, C->D, ... J->K); whatever used to sit in the row's last
# used column falls off the fixed A:K used range. Column A (labels) and
# row 1 (the Q0..Q9 headers) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New value to drop into column B for each data row.
$newValues = @{
    2  = 2.249802839611392
    3  = 7.469150330857293
    4  = -18.17126180013747
    5  = 8.13967154697915
    6  = 2.147570471799392
    7  = -3.44002998652333
    8  = -0.06117417841420103
    9  = 2.041455937156254
    10 = -0.5825945370336409
    11 = 0.09567504080935779
    12 = -0.2604190369987228
    13 = 0.8354549961584912
    14 = -0.1000793599026215
    15 = -0.3537865060796963
    16 = 0.1481773904324453
    17 = 0.157445989004155
    18 = -0.5006594565260708
    19 = 0.2803578805354692
    20 = -0.1719748578450117
    21 = 0.3058625397463315
    22 = -0.6123299526872862
    23 = 0.6883713851991116
    24 = -0.2766911554241067
}

# Last used column (1-based, A=1) of the pre-edit data in each row - i.e.
# how far right the existing B:K values extend before the shift.
$lastCol = @{
    2 = 11; 3 = 11; 4 = 11; 5 = 11; 6 = 11; 7 = 11; 8 = 11; 9 = 11; 10 = 11
    11 = 11; 12 = 11; 13 = 11; 14 = 11
    15 = 10; 16 = 9; 17 = 8; 18 = 7; 19 = 6; 20 = 5; 21 = 4; 22 = 3; 23 = 2; 24 = 1
}

for ($row = 2; $row -le 24; $row++) {
    $last = $lastCol[$row]
    # Shift existing values one column to the right, working from the
    # rightmost populated column back down to column B (col 2) so we never
    # overwrite a value before it has been read.
    for ($col = $last; $col -ge 2; $col--) {
        $val = $ws.Cells.Item($row, $col).Value2
        $ws.Cells.Item($row, $col + 1).Value2 = $val
    }
    # Whatever previously sat in the row's last used column has now been
    # duplicated one column further right (at $last + 1); that overflow
    # falls off the sheet's fixed A:K used range, so clear the spill-over
    # cell rather than leaving a stray value behind.
    $ws.Cells.Item($row, $last + 1).ClearContents()

    # Write the newly computed value into the now-vacated column B.
    $ws.Cells.Item($row, 2).Value2 = $newValues[$row]
}
